# week 2 spring 24
# Update the "Nine" matchup-averages sheet with the latest computed
# averages after adding week 2 spring '24 results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$updates = @{
    "D2" = 11.28
    "D3" = 10.25
    "E3" = 10.63
    "F3" = 10.07
    "B4" = 8.720000000000001
    "C4" = 9.75
    "E4" = 10.58
    "F4" = 10.13
    "H4" = 10.5
    "C5" = 9.369999999999999
    "D5" = 9.42
    "F5" = 10.3
    "G5" = 9.609999999999999
    "C6" = 9.93
    "D6" = 9.869999999999999
    "E6" = 9.699999999999999
    "G6" = 10.37
    "I6" = 8.94
    "J6" = 7.71
    "E7" = 10.39
    "F7" = 9.630000000000001
    "H7" = 9.56
    "I7" = 7.17
    "D8" = 9.5
    "G8" = 10.44
    "F9" = 11.06
    "G9" = 12.83
    "F10" = 12.29
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
